$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Append new variable rows (9-15) below the existing variable table ---
$labels = @("emissivity", "rc", "mu_0", "Cp_0", "Pr_0", "wall_thick", "inregen_wall_thick")
# High-precision decimal literals chosen so the stored IEEE-754 double exactly
# matches the intended values (0.06, 0.003, 6.848e-5, 1.9, 0.5758, 3e-4, 1e-3).
$values = @(
    0.059999999999999997779553950750,
    0.003000000000000000062450045135,
    0.000068480000000000008420798780,
    1.899999999999999911182158029987,
    0.575799999999999978506082243257,
    0.000299999999999999973718939339,
    0.001000000000000000020816681712
)

for ($i = 0; $i -lt $labels.Length; $i++) {
    $row = 9 + $i
    $ws.Range("A$row").Value = $labels[$i]
    $ws.Range("B$row").Value = $values[$i]
}

# --- Style the column H marker cells for rows 8-12 ---
# Reuse the existing Consolas-based style already present on N3, then recolor
# it to the new reddish accent tone used for these new marker cells.
$srcFormat = $ws.Range("N3")
$srcFormat.Copy()
$markerRows = 8, 9, 10, 11, 12
foreach ($r in $markerRows) {
    $cell = $ws.Range("H$r")
    $cell.PasteSpecial(-4122)
    $cell.Font.Color = 7695584
}
$excel.CutCopyMode = 0

# --- Column widths for B and L to fit the new content ---
$ws.Columns.Item(2).ColumnWidth = 10.14
$ws.Columns.Item(12).ColumnWidth = 10.14

# --- Update the active selection to match the saved workbook state ---
$ws.Range("B16").Select()
